$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 64: correct the debited amount (was a placeholder of 15600) ---
$ws.Range("B64").Value = 67600
# Row 64 is no longer the last "Ordered Amount" entry, so drop the italic marker
$ws.Range("D64").Font.Italic = $false

# --- Row 65: new ledger entry for 17-Feb-2020 ---
$ws.Range("A65").Value = 43878
$ws.Range("A65").NumberFormat = "[$-409]d\-mmm\-yyyy;@"
$ws.Range("A65").HorizontalAlignment = -4108
$ws.Range("B65").Value = 38480
$ws.Range("D65").Value = "Ordered Amount"
$ws.Range("D65").HorizontalAlignment = -4108
$ws.Range("D65").Font.Italic = $false
$ws.Range("E65").Formula = "=IF(A65=`"`",`"`",SUM(E64-B65+C65))"

# --- Row 66: new ledger entry for 18-Feb-2020 (now the latest / last row) ---
$ws.Range("A66").Value = 43879
$ws.Range("A66").NumberFormat = "[$-409]d\-mmm\-yyyy;@"
$ws.Range("A66").HorizontalAlignment = -4108
$ws.Range("B66").Value = 11440
$ws.Range("D66").Value = "Ordered Amount"
$ws.Range("D66").HorizontalAlignment = -4108
$ws.Range("D66").Font.Italic = $true
$ws.Range("E66").Formula = "=IF(A66=`"`",`"`",SUM(E65-B66+C66))"

# --- Scroll the frozen pane down and move the selection to the new last entry ---
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("E66").Select()
